$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.239.37'
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").Value = '2.281.89'
$ws.Range("E3").Value = '  +3.31%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.43'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.20'
$ws.Range("E6").Value = '  +6.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.590'
$ws.Range("E7").Value = '  +1.89%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.573'
$ws.Range("E9").Value = '  +3.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.06'
$ws.Range("E10").Value = '  +7.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("E11").Value = '  +2.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.92'
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("D14").Value = '2.626.71'
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.883'
$ws.Range("E15").Value = '  +2.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.65'
$ws.Range("E16").Value = '  +4.22%  '
$ws.Range("D17").Value = '2.284.57'
$ws.Range("E17").Value = '  +3.73%  '
$ws.Range("D18").Value = '44.140.02'
$ws.Range("E18").Value = '  +3.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.16'
$ws.Range("E19").Value = '  -4.13%  '
$ws.Range("D20").Value = '0.0₃0999'
$ws.Range("E20").Value = '  +4.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.60'
$ws.Range("E21").Value = '  +3.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.44'
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.23'
$ws.Range("E23").Value = '  +3.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '238.22'
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("E25").Value = '  +4.25%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("E27").Value = '  +2.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.29'
$ws.Range("E28").Value = '  +15.80%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.53'
$ws.Range("E30").Value = '  +5.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.68'
$ws.Range("E31").Value = '  +5.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.56'
$ws.Range("E32").Value = '  +1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0884'
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.09'
$ws.Range("E35").Value = '  +5.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.29'
$ws.Range("E36").Value = '  +4.03%  '
$ws.Range("E37").Value = '  +10.81%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.53'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.94'
$ws.Range("E40").Value = '  +6.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.68'
$ws.Range("E41").Value = '  +29.48%  '
$ws.Range("E42").Value = '  +0.86%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '1.773.53'
$ws.Range("E44").Value = '  -4.88%  '
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '85.66'
$ws.Range("E46").Value = '  -3.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.41'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.96'
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.57'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.87'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '104.88'
$ws.Range("E51").Value = '  +4.19%  '
